$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rule "R40" (row 11) is renamed to "1": the rule-name cell B11 changes its
# text from "R40" to "1". A plain "$ws.Range('B11').Value = '1'" would let
# Excel auto-detect the numeric-looking text as a *number*, which also
# drags in a new number-format style. To keep B11 a genuine text value (and
# keep its existing style/format untouched, exactly like the source edit),
# build the text "1" in a scratch cell via a formula (so it is never
# auto-converted to a number), copy it, and paste only the value into B11.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = "=""1"""
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

